$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.977669497583861
$ws.Range("J2").Value = 0.977669497583861
$ws.Range("M2").Value = 144.2773156666667
$ws.Range("N2").Value = 432.831947
$ws.Range("O2").Value = 0.8052971554812057
$ws.Range("P2").Value = 0.8052971554812056
$ws.Range("Q2").Value = 1084.214450385288
$ws.Range("R2").Value = 9757.930053467595
$ws.Range("S2").Value = 0.7873144654050228
$ws.Range("T2").Value = 0.7873144654050227

# Row 3
$ws.Range("I3").Value = 0.977669497583861
$ws.Range("J3").Value = 0.977669497583861
$ws.Range("O3").Value = 0.0082793637854752
$ws.Range("P3").Value = 0.008279363785475198
$ws.Range("S3").Value = 0.008094481432459553
$ws.Range("T3").Value = 0.008094481432459551

# Row 4
$ws.Range("I4").Value = 0.977669497583861
$ws.Range("J4").Value = 0.977669497583861
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5156256666666666
$ws.Range("N4").Value = 1.546877
$ws.Range("O4").Value = 0.002878012255364554
$ws.Range("P4").Value = 0.002878012255364553
$ws.Range("Q4").Value = 3.874821181738333
$ws.Range("R4").Value = 34.873390635645
$ws.Range("S4").Value = 0.002813744795742458
$ws.Range("T4").Value = 0.002813744795742458

# Row 5
$ws.Range("I5").Value = 0.977669497583861
$ws.Range("J5").Value = 0.977669497583861
$ws.Range("M5").Value = 30.00245966666667
$ws.Range("N5").Value = 90.007379
$ws.Range("O5").Value = 0.1674614981250883
$ws.Range("P5").Value = 0.1674614981250883
$ws.Range("Q5").Value = 225.4623338907683
$ws.Range("R5").Value = 2029.161005016915
$ws.Range("S5").Value = 0.1637219987365957
$ws.Range("T5").Value = 0.1637219987365957

# Row 6
$ws.Range("I6").Value = 0.977669497583861
$ws.Range("J6").Value = 0.977669497583861
$ws.Range("M6").Value = 2.189762333333333
$ws.Range("N6").Value = 6.569287
$ws.Range("O6").Value = 0.01222236059816459
$ws.Range("P6").Value = 0.01222236059816459
$ws.Range("Q6").Value = 16.45561503372167
$ws.Range("R6").Value = 148.100535303495
$ws.Range("S6").Value = 0.01194942914529635
$ws.Range("T6").Value = 0.01194942914529635

# Row 7
$ws.Range("I7").Value = 0.977669497583861
$ws.Range("J7").Value = 0.977669497583861
$ws.Range("M7").Value = 0.6918473333333334
$ws.Range("N7").Value = 2.075542
$ws.Range("O7").Value = 0.0038616097547018
$ws.Range("P7").Value = 0.0038616097547018
$ws.Range("Q7").Value = 5.199090881296667
$ws.Range("R7").Value = 46.79181793167
$ws.Range("S7").Value = 0.003775378068744246
$ws.Range("T7").Value = 0.003775378068744245

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.171642
$ws.Range("H8").Value = 0.514926
$ws.Range("I8").Value = 0.02233050241613897
$ws.Range("J8").Value = 0.02233050241613898
$ws.Range("M8").Value = 144.2773156666667
$ws.Range("N8").Value = 432.831947
$ws.Range("O8").Value = 0.8052971554812057
$ws.Range("P8").Value = 0.8052971554812056
$ws.Range("Q8").Value = 24.764047015658
$ws.Range("R8").Value = 222.876423140922
$ws.Range("S8").Value = 0.01798269007618291
$ws.Range("T8").Value = 0.01798269007618291

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.171642
$ws.Range("H9").Value = 0.514926
$ws.Range("I9").Value = 0.02233050241613897
$ws.Range("J9").Value = 0.02233050241613898
$ws.Range("O9").Value = 0.0082793637854752
$ws.Range("P9").Value = 0.008279363785475198
$ws.Range("Q9").Value = 0.254602357214
$ws.Range("R9").Value = 2.291421214926
$ws.Range("S9").Value = 0.0001848823530156475
$ws.Range("T9").Value = 0.0001848823530156475

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.171642
$ws.Range("H10").Value = 0.514926
$ws.Range("I10").Value = 0.02233050241613897
$ws.Range("J10").Value = 0.02233050241613898
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5156256666666666
$ws.Range("N10").Value = 1.546877
$ws.Range("O10").Value = 0.002878012255364554
$ws.Range("P10").Value = 0.002878012255364553
$ws.Range("Q10").Value = 0.08850302067799999
$ws.Range("R10").Value = 0.796527186102
$ws.Range("S10").Value = 0.00006426745962209574
$ws.Range("T10").Value = 0.00006426745962209574

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.171642
$ws.Range("H11").Value = 0.514926
$ws.Range("I11").Value = 0.02233050241613897
$ws.Range("J11").Value = 0.02233050241613898
$ws.Range("M11").Value = 30.00245966666667
$ws.Range("N11").Value = 90.007379
$ws.Range("O11").Value = 0.1674614981250883
$ws.Range("P11").Value = 0.1674614981250883
$ws.Range("Q11").Value = 5.149682182106
$ws.Range("R11").Value = 46.347139638954
$ws.Range("S11").Value = 0.003739499388492536
$ws.Range("T11").Value = 0.003739499388492536

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.171642
$ws.Range("H12").Value = 0.514926
$ws.Range("I12").Value = 0.02233050241613897
$ws.Range("J12").Value = 0.02233050241613898
$ws.Range("M12").Value = 2.189762333333333
$ws.Range("N12").Value = 6.569287
$ws.Range("O12").Value = 0.01222236059816459
$ws.Range("P12").Value = 0.01222236059816459
$ws.Range("Q12").Value = 0.375855186418
$ws.Range("R12").Value = 3.382696677762
$ws.Range("S12").Value = 0.0002729314528682361
$ws.Range("T12").Value = 0.0002729314528682361

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.171642
$ws.Range("H13").Value = 0.514926
$ws.Range("I13").Value = 0.02233050241613897
$ws.Range("J13").Value = 0.02233050241613898
$ws.Range("M13").Value = 0.6918473333333334
$ws.Range("N13").Value = 2.075542
$ws.Range("O13").Value = 0.0038616097547018
$ws.Range("P13").Value = 0.0038616097547018
$ws.Range("Q13").Value = 0.118750059988
$ws.Range("R13").Value = 1.068750539892
$ws.Range("S13").Value = 0.00008623168595755438
$ws.Range("T13").Value = 0.00008623168595755437
